# actor 테이블에 attackRange, ultimateRange 추가
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActorTable")

# New header values
$ws.Range("H1").Value = "attackRange|Float"
$ws.Range("I1").Value = "ultimateRange|Float"

# attackRange values for rows 2..18
$attackRange = @(0, 4, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
# ultimateRange values for rows 2..18 (all zero)
$ultimateRange = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt 17; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $attackRange[$i]
    $ws.Cells.Item($row, 9).Value = $ultimateRange[$i]
}
